# Revert "feat(dialog): update CN data and dialogue Excel files"
# - Change B3 version string back from "EA 23.252" to "EA 23.200"
# - Remove the "flare_" row (originally row 11)
# - Remove the "SpMoonArrow" row (originally row 17)
# - Remove the trailing "SpReturn" / "SpEvac" rows (originally rows 37-38)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows bottom-to-top so earlier row numbers stay valid.
$ws.Rows(38).Delete()
$ws.Rows(37).Delete()
$ws.Rows(17).Delete()
$ws.Rows(11).Delete()

# Restore the old version string for SpHolyVeil.
$ws.Range("B3").Value = "EA 23.200"
